$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.084564685821533
$ws.Range("B1").Value = 6.216755390167236
$ws.Range("C1").Value = 6.211811065673828
$ws.Range("D1").Value = 6.584945678710938
$ws.Range("E1").Value = 5.333799362182617
